# Prevent storekeepers from directly issuing materials without approval.
# Simulates a "White Tiles" issuance of 4 boxes going "To tiler":
#   - Current Stock: decrements White Tiles quantity from 5 -> 1,
#     updates the Last Updated timestamp and Updated By to "storesperson".
#   - Issuance Log: appends the corresponding issuance row.

$wb = $excel.ActiveWorkbook

$stockSheet = $wb.Worksheets.Item("Current Stock")
$logSheet   = $wb.Worksheets.Item("Issuance Log")

$timestamp = "2025-07-10 12:19:24"

# --- Update "Current Stock" row for White Tiles (row 2) ---
$stockSheet.Range("B2").Value = 1
$stockSheet.Range("D2").Value = $timestamp
$stockSheet.Range("E2").Value = "storesperson"

# --- Append new row to "Issuance Log" ---
$logSheet.Range("A2").Value = $timestamp
$logSheet.Range("B2").Value = "White Tiles"
$logSheet.Range("C2").Value = 4
$logSheet.Range("D2").Value = "boxes"
$logSheet.Range("E2").Value = "storesperson"
$logSheet.Range("F2").Value = "To tiler"
$logSheet.Range("G2").Value = 1
